$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"25.8200000000006"
$ws.Range("H2").Value = [double]"1.049284770182624e-09"
$ws.Range("I2").Value = [double]"1.049284770182624e-09"
$ws.Range("L2").Value = [double]"61.1627775250962"
$ws.Range("M2").Value = "[44.851213955111746, 77.47434109508065]"
$ws.Range("N2").Value = [double]"1.557985296329889e-09"
$ws.Range("O2").Value = [double]"1.557985296329889e-09"
$ws.Range("P2").Value = [double]"1.540921321580579"
$ws.Range("Q2").Value = "[1.2138686329185786, 1.8679740102425795]"
$ws.Range("R2").Value = [double]"2.615019312202094e-12"
$ws.Range("S2").Value = [double]"2.615019312202094e-12"
$ws.Range("T2").Value = [double]"67.11900601826797"
$ws.Range("U2").Value = "[56.16525113780611, 78.07276089872983]"
$ws.Range("V2").Value = [double]"4.440892098500626e-16"
$ws.Range("W2").Value = [double]"4.440892098500626e-16"
$ws.Range("X2").Value = [double]"19.48776776776822"
$ws.Range("Y2").Value = [double]"18.1437837837842"
$ws.Range("Z2").Value = [double]"20.83175175175223"

# Row 3
$ws.Range("F3").Value = [double]"25.8200000000006"
$ws.Range("H3").Value = [double]"4.113947860151512e-08"
$ws.Range("I3").Value = [double]"4.113947860151512e-08"
$ws.Range("L3").Value = [double]"54.9843348609207"
$ws.Range("M3").Value = "[36.75985316153579, 73.2088165603056]"
$ws.Range("N3").Value = [double]"2.402335661511046e-07"
$ws.Range("O3").Value = [double]"2.402335661511046e-07"
$ws.Range("P3").Value = [double]"1.490605523324887"
$ws.Range("Q3").Value = "[1.0880791372793475, 1.8931319093704264]"
$ws.Range("R3").Value = [double]"2.140251975646379e-09"
$ws.Range("S3").Value = [double]"2.140251975646379e-09"
$ws.Range("T3").Value = [double]"58.36480980165975"
$ws.Range("U3").Value = "[47.034785305407596, 69.69483429791191]"
$ws.Range("V3").Value = [double]"1.618705169903478e-13"
$ws.Range("W3").Value = [double]"1.618705169903478e-13"
$ws.Range("X3").Value = [double]"19.69453453453499"
$ws.Range("Y3").Value = [double]"18.04040040040081"
$ws.Range("Z3").Value = [double]"21.34866866866916"

# Row 4
$ws.Range("F4").Value = [double]"25.8200000000006"
$ws.Range("H4").Value = [double]"2.531863607657669e-12"
$ws.Range("I4").Value = [double]"2.531863607657669e-12"
$ws.Range("L4").Value = [double]"64.35445729283565"
$ws.Range("M4").Value = "[47.14507931080309, 81.56383527486821]"
$ws.Range("N4").Value = [double]"1.669731242159855e-09"
$ws.Range("O4").Value = [double]"1.669731242159855e-09"
$ws.Range("P4").Value = [double]"0.6100790538502707"
$ws.Range("Q4").Value = "[0.35850006257180755, 0.8616580451287339]"
$ws.Range("R4").Value = [double]"1.351329201737528e-05"
$ws.Range("S4").Value = [double]"1.351329201737528e-05"
$ws.Range("T4").Value = [double]"58.59161875798431"
$ws.Range("U4").Value = "[49.80039509722473, 67.38284241874388]"
$ws.Range("V4").Value = [double]"0"
$ws.Range("W4").Value = [double]"0"
$ws.Range("X4").Value = [double]"23.31295295295349"
$ws.Range("Y4").Value = [double]"22.27911911911962"
$ws.Range("Z4").Value = [double]"24.34678678678735"

# Row 5
$ws.Range("B5").Value = [double]"1"
$ws.Range("F5").Value = [double]"25.8200000000006"
$ws.Range("H5").Value = [double]"1.197147166653867e-05"
$ws.Range("I5").Value = [double]"1.197147166653867e-05"
$ws.Range("L5").Value = [double]"45.06753531122817"
$ws.Range("M5").Value = "[23.791366086891472, 66.34370453556487]"
$ws.Range("N5").Value = [double]"0.0001008468872447654"
$ws.Range("O5").Value = [double]"0.0001008468872447654"
$ws.Range("P5").Value = [double]"0.4968685077749626"
$ws.Range("Q5").Value = "[0.006289474781961957, 0.9874475407679633]"
$ws.Range("R5").Value = [double]"0.04725353075123873"
$ws.Range("S5").Value = [double]"0.04725353075123873"
$ws.Range("T5").Value = [double]"49.10191928637165"
$ws.Range("U5").Value = "[37.6835725560052, 60.52026601673809]"
$ws.Range("V5").Value = [double]"3.836886364183556e-11"
$ws.Range("W5").Value = [double]"3.836886364183556e-11"
$ws.Range("X5").Value = [double]"23.77817817817873"
$ws.Range("Y5").Value = [double]"21.7622022022027"
$ws.Range("Z5").Value = [double]"25.79415415415475"

# Row 6
$ws.Range("F6").Value = [double]"25.8200000000006"
$ws.Range("H6").Value = [double]"3.331002584872067e-09"
$ws.Range("I6").Value = [double]"3.331002584872067e-09"
$ws.Range("L6").Value = [double]"61.10164790641871"
$ws.Range("M6").Value = "[44.53093650583533, 77.67235930700208]"
$ws.Range("N6").Value = [double]"2.384159092372329e-09"
$ws.Range("O6").Value = [double]"2.384159092372329e-09"
$ws.Range("P6").Value = [double]"0.03144737390980801"
$ws.Range("Q6").Value = "[-0.2704474156243464, 0.3333421634439624]"
$ws.Range("R6").Value = [double]"0.8347688676203104"
$ws.Range("S6").Value = [double]"0.8347688676203104"
$ws.Range("T6").Value = [double]"62.75166979219343"
$ws.Range("U6").Value = "[52.075306485870286, 73.42803309851658]"
$ws.Range("V6").Value = [double]"1.998401444325282e-15"
$ws.Range("W6").Value = [double]"1.998401444325282e-15"
$ws.Range("X6").Value = [double]"25.69077077077137"
$ws.Range("Y6").Value = [double]"24.45017017017074"
$ws.Range("Z6").Value = [double]"26.93137137137199"

# Row 7
$ws.Range("F7").Value = [double]"25.8200000000006"
$ws.Range("H7").Value = [double]"1.626779466690209e-08"
$ws.Range("I7").Value = [double]"1.626779466690209e-08"
$ws.Range("L7").Value = [double]"56.94587460421737"
$ws.Range("M7").Value = "[36.27952010421495, 77.6122291042198]"
$ws.Range("N7").Value = [double]"1.447086214190207e-06"
$ws.Range("O7").Value = [double]"1.447086214190207e-06"
$ws.Range("P7").Value = [double]"-0.6415264277600778"
$ws.Range("Q7").Value = "[-1.0063159651138474, -0.27673689040630833]"
$ws.Range("R7").Value = [double]"0.0009376536375877986"
$ws.Range("S7").Value = [double]"0.0009376536375877986"
$ws.Range("T7").Value = [double]"66.52334666875868"
$ws.Range("U7").Value = "[55.41361564958129, 77.63307768793607]"
$ws.Range("V7").Value = [double]"1.110223024625157e-15"
$ws.Range("W7").Value = [double]"1.110223024625157e-15"
$ws.Range("X7").Value = [double]"2.636276276276337"
$ws.Range("Y7").Value = [double]"1.137217217217245"
$ws.Range("Z7").Value = [double]"4.13533533533543"

# Row 8
$ws.Range("F8").Value = [double]"25.8200000000006"
$ws.Range("H8").Value = [double]"2.526061582130978e-10"
$ws.Range("I8").Value = [double]"2.526061582130978e-10"
$ws.Range("L8").Value = [double]"52.13223070074569"
$ws.Range("M8").Value = "[35.03731329129526, 69.22714811019613]"
$ws.Range("N8").Value = [double]"1.920086545315058e-07"
$ws.Range("O8").Value = [double]"1.920086545315058e-07"
$ws.Range("P8").Value = [double]"-0.9937370155499243"
$ws.Range("Q8").Value = "[-1.3207897042119248, -0.6666843268879239]"
$ws.Range("R8").Value = [double]"2.072936433616945e-07"
$ws.Range("S8").Value = [double]"2.072936433616945e-07"
$ws.Range("T8").Value = [double]"50.11537211740885"
$ws.Range("U8").Value = "[41.24844618378312, 58.98229805103457]"
$ws.Range("V8").Value = [double]"7.771561172376096e-15"
$ws.Range("W8").Value = [double]"7.771561172376096e-15"
$ws.Range("X8").Value = [double]"4.08364364364374"
$ws.Range("Y8").Value = [double]"2.739659659659725"
$ws.Range("Z8").Value = [double]"5.427627627627754"

# Row 9
$ws.Range("F9").Value = [double]"24.80000000000044"
$ws.Range("H9").Value = [double]"9.698369238808624e-07"
$ws.Range("I9").Value = [double]"9.698369238808624e-07"
$ws.Range("L9").Value = [double]"47.93625960279892"
$ws.Range("M9").Value = "[28.99608756737635, 66.87643163822149]"
$ws.Range("N9").Value = [double]"6.643167265174554e-06"
$ws.Range("O9").Value = [double]"6.643167265174554e-06"
$ws.Range("P9").Value = [double]"-1.333368653775848"
$ws.Range("Q9").Value = "[-1.7987897876410024, -0.8679475199106932]"
$ws.Range("R9").Value = [double]"6.842369264514758e-07"
$ws.Range("S9").Value = [double]"6.842369264514758e-07"
$ws.Range("T9").Value = [double]"50.39389331197469"
$ws.Range("U9").Value = "[39.1655766483136, 61.62220997563578]"
$ws.Range("V9").Value = [double]"1.11441966765824e-11"
$ws.Range("W9").Value = [double]"1.11441966765824e-11"
$ws.Range("X9").Value = [double]"5.262862862862956"
$ws.Range("Y9").Value = [double]"3.425825825825886"
$ws.Range("Z9").Value = [double]"7.099899899900026"

# Row 10
$ws.Range("F10").Value = [double]"24.80000000000044"
$ws.Range("H10").Value = [double]"2.676225741016225e-08"
$ws.Range("I10").Value = [double]"2.676225741016225e-08"
$ws.Range("L10").Value = [double]"55.78083033352755"
$ws.Range("M10").Value = "[38.92215239352072, 72.63950827353439]"
$ws.Range("N10").Value = [double]"3.213507882016131e-08"
$ws.Range("O10").Value = [double]"3.213507882016131e-08"
$ws.Range("P10").Value = [double]"-1.572368695490387"
$ws.Range("Q10").Value = "[-1.9497371824080805, -1.1950002085726936]"
$ws.Range("R10").Value = [double]"9.32867116887337e-11"
$ws.Range("S10").Value = [double]"9.32867116887337e-11"
$ws.Range("T10").Value = [double]"52.11831158006915"
$ws.Range("U10").Value = "[40.99032502244504, 63.246298137693266]"
$ws.Range("V10").Value = [double]"3.132161197072492e-12"
$ws.Range("W10").Value = [double]"3.132161197072492e-12"
$ws.Range("X10").Value = [double]"6.206206206206318"
$ws.Range("Y10").Value = [double]"4.7167167167168"
$ws.Range("Z10").Value = [double]"7.695695695695836"

# Row 11
$ws.Range("F11").Value = [double]"24.80000000000044"
$ws.Range("H11").Value = [double]"1.002289251594846e-09"
$ws.Range("I11").Value = [double]"1.002289251594846e-09"
$ws.Range("L11").Value = [double]"58.21721892863206"
$ws.Range("M11").Value = "[42.680441165953525, 73.7539966913106]"
$ws.Range("N11").Value = [double]"1.585836795214846e-09"
$ws.Range("O11").Value = [double]"1.585836795214846e-09"
$ws.Range("P11").Value = [double]"-1.509473947670771"
$ws.Range("Q11").Value = "[-1.8239476867688484, -1.1950002085726936]"
$ws.Range("R11").Value = [double]"1.482813871689359e-12"
$ws.Range("S11").Value = [double]"1.482813871689359e-12"
$ws.Range("T11").Value = [double]"64.09515611942719"
$ws.Range("U11").Value = "[53.871008394049, 74.31930384480538]"
$ws.Range("V11").Value = [double]"2.220446049250313e-16"
$ws.Range("W11").Value = [double]"2.220446049250313e-16"
$ws.Range("X11").Value = [double]"5.957957957958062"
$ws.Range("Y11").Value = [double]"4.716716716716799"
$ws.Range("Z11").Value = [double]"7.199199199199325"
